$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 809.7222
$ws.Range("I19").Value = 1030.25
$ws.Range("J19").Value = 746.7143
$ws.Range("K19").Value = 1030.25
$ws.Range("L19").Value = 746.7143
$ws.Range("M19").Value = -855.25
$ws.Range("N19").Value = -1096.7143
$ws.Range("H110").Value = 25080.455
$ws.Range("J110").Value = 25080.455
$ws.Range("L110").Value = 25080.455
$ws.Range("N110").Value = -33260.455
$ws.Range("H123").Value = 48188.57
$ws.Range("J123").Value = 48188.57
$ws.Range("L123").Value = 48188.57
$ws.Range("N123").Value = -57988.57
$ws.Range("H128").Value = 500024000
$ws.Range("J128").Value = 500024000
$ws.Range("L128").Value = 500024000
$ws.Range("N128").Value = -500033960
$ws.Range("H132").Value = 3859.9656
$ws.Range("I132").Value = 3014.4531
$ws.Range("J132").Value = 6212.696
$ws.Range("K132").Value = 9043.3593
$ws.Range("L132").Value = 18638.088
$ws.Range("M132").Value = -6513.3593
$ws.Range("N132").Value = -23698.088
$ws.Range("H135").Value = 671.10205
$ws.Range("I135").Value = 372.97562
$ws.Range("J135").Value = 2199
$ws.Range("K135").Value = 3356.78058
$ws.Range("L135").Value = 19791
$ws.Range("M135").Value = -821.7805800000001
$ws.Range("N135").Value = -24861
$ws.Range("H138").Value = 1734.7727
$ws.Range("I138").Value = 985.75
$ws.Range("J138").Value = 3732.1667
$ws.Range("K138").Value = 2957.25
$ws.Range("L138").Value = 11196.5001
$ws.Range("M138").Value = 2182.75
$ws.Range("N138").Value = -21476.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 1897953.9
$ws.Range("I132").Value = 2805.2058
$ws.Range("J132").Value = 3911549.2
$ws.Range("K132").Value = 8415.617400000001
$ws.Range("L132").Value = 11734647.6
$ws.Range("M132").Value = -5885.617400000001
$ws.Range("N132").Value = -11739707.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 271.85715
$ws.Range("I80").Value = 298.0625
$ws.Range("J80").Value = 236.91667
$ws.Range("K80").Value = 298.0625
$ws.Range("L80").Value = 236.91667
$ws.Range("M80").Value = 699.9375
$ws.Range("N80").Value = -2232.91667
$ws.Range("H83").Value = 271.85715
$ws.Range("I83").Value = 298.0625
$ws.Range("J83").Value = 236.91667
$ws.Range("K83").Value = 1490.3125
$ws.Range("L83").Value = 1184.58335
$ws.Range("M83").Value = 3501.6875
$ws.Range("N83").Value = -11168.58335
$ws.Range("H94").Value = 1058.1111
$ws.Range("I94").Value = 1003.06665
$ws.Range("J94").Value = 1333.3334
$ws.Range("K94").Value = 1003.06665
$ws.Range("L94").Value = 1333.3334
$ws.Range("M94").Value = -552.06665
$ws.Range("N94").Value = -2235.3334
$ws.Range("H134").Value = 3933.0923
$ws.Range("I134").Value = 1506.5526
$ws.Range("J134").Value = 7348.222
$ws.Range("K134").Value = 4519.6578
$ws.Range("L134").Value = 22044.666
$ws.Range("M134").Value = -1984.6578
$ws.Range("N134").Value = -27114.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2312.0227
$ws.Range("I132").Value = 1425.96
$ws.Range("J132").Value = 3477.8948
$ws.Range("K132").Value = 4277.88
$ws.Range("L132").Value = 10433.6844
$ws.Range("M132").Value = -1747.88
$ws.Range("N132").Value = -15493.6844
$ws.Range("H134").Value = 1061.921
$ws.Range("I134").Value = 638.24
$ws.Range("J134").Value = 1876.6923
$ws.Range("K134").Value = 1914.72
$ws.Range("L134").Value = 5630.0769
$ws.Range("M134").Value = 620.28
$ws.Range("N134").Value = -10700.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 500785.03
$ws.Range("J131").Value = 834381.75
$ws.Range("L131").Value = 2503145.25
$ws.Range("N131").Value = -2513225.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2552505.2
$ws.Range("J21").Value = 70007
$ws.Range("L21").Value = 70007
$ws.Range("N21").Value = -70353
$ws.Range("H30").Value = 2552505.2
$ws.Range("J30").Value = 70007
$ws.Range("L30").Value = 70007
$ws.Range("N30").Value = -70217
$ws.Range("H132").Value = 564885.25
$ws.Range("I132").Value = 851534.2
$ws.Range("J132").Value = 3053.28
$ws.Range("K132").Value = 2554602.6
$ws.Range("L132").Value = 9159.84
$ws.Range("M132").Value = -2552072.6
$ws.Range("N132").Value = -14219.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1401
$ws.Range("I93").Value = 1471.2106
$ws.Range("J93").Value = 1067.5
$ws.Range("K93").Value = 1471.2106
$ws.Range("L93").Value = 1067.5
$ws.Range("M93").Value = -223.2106000000001
$ws.Range("N93").Value = -3563.5
$ws.Range("H122").Value = 3497.6
$ws.Range("I122").Value = 2458.4
$ws.Range("J122").Value = 5576
$ws.Range("K122").Value = 7375.200000000001
$ws.Range("L122").Value = 16728
$ws.Range("M122").Value = -4925.200000000001
$ws.Range("N122").Value = -21628
$ws.Range("H128").Value = 333359970
$ws.Range("J128").Value = 333359970
$ws.Range("L128").Value = 333359970
$ws.Range("N128").Value = -333369930
$ws.Range("H132").Value = 20002872
$ws.Range("I132").Value = 43482316
$ws.Range("J132").Value = 1866.0741
$ws.Range("K132").Value = 130446948
$ws.Range("L132").Value = 5598.2223
$ws.Range("M132").Value = -130444418
$ws.Range("N132").Value = -10658.2223
$ws.Range("H136").Value = 8475774
$ws.Range("I136").Value = 11111701
$ws.Range("J136").Value = 3150
$ws.Range("K136").Value = 33335103
$ws.Range("L136").Value = 9450
$ws.Range("M136").Value = -33332553
$ws.Range("N136").Value = -14550

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1346.8182
$ws.Range("I81").Value = 729.375
$ws.Range("K81").Value = 1458.75
$ws.Range("M81").Value = -397.75
$ws.Range("H84").Value = 1346.8182
$ws.Range("I84").Value = 729.375
$ws.Range("K84").Value = 7293.75
$ws.Range("M84").Value = -1989.75
$ws.Range("H132").Value = 1547.3422
$ws.Range("I132").Value = 1393.8334
$ws.Range("J132").Value = 1736.9706
$ws.Range("K132").Value = 4181.5002
$ws.Range("L132").Value = 5210.9118
$ws.Range("M132").Value = -1651.5002
$ws.Range("N132").Value = -10270.9118
$ws.Range("H136").Value = 4391294
$ws.Range("I136").Value = 5440699
$ws.Range("J136").Value = 2872.7273
$ws.Range("K136").Value = 16322097
$ws.Range("L136").Value = 8618.1819
$ws.Range("M136").Value = -16319547
$ws.Range("N136").Value = -13718.1819
